$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 currently holds the three most recent "blog" tiles referencing
# ser: 160 (H10), ser: 161 (D10) and ser: 162 (B10). A new blog post
# (ser: 163) was added, so the window slides forward by one: the oldest
# entry (ser: 160) drops off and each remaining entry shifts into the
# next slot, with the new post taking the first slot (B10).

$ws.Range("H10").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 161"
$ws.Range("D10").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 162"
$ws.Range("B10").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 163"
